$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "MONTHLY MEAL DATA"

# Update row 2 (employee John Doe, now first record)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "John Doe"
$ws.Range("C2").Value = "Normal"
$ws.Range("D2").Value = "3:49 PM"
$ws.Range("E2").Value = "15 July 2024"
$ws.Range("F2").Value = "Karachi"

# Update row 3 (same employee, second meal entry)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "John Doe"
$ws.Range("C3").Value = "Normal"
$ws.Range("D3").Value = "4:02 PM"
$ws.Range("E3").Value = "16 July 2024"
$ws.Range("F3").Value = "Karachi"

# Delete rows 4 and 5 (old extra records no longer present)
$ws.Range("A4:F5").Delete()
